{"js": "// Expand each single weekly placeholder in the Rencana Pembelajaran table\n// into 16 per-row placeholders (one per week). Each placeholder keeps\n// its own line inside the very same run (a literal line-feed inside the\n// <w:t>), matching the template's new per-row placeholder scheme, while\n// the run's original formatting (rFonts/sz/szCs/lang) is preserved\n// unchanged.\n\n// old placeholder name -> new per-week prefix\nconst map = [\n  [\"MINGGU_KE\", \"MINGGU\"],\n  [\"SUB_CPMK_MINGGUAN\", \"SUB_CPMK\"],\n  [\"INDIKATOR_MINGGUAN\", \"INDIKATOR\"],\n  [\"TOPIK_MINGGUAN\", \"TOPIK\"],\n  [\"METODE_PEMBELAJARAN\", \"METODE\"],\n  [\"ESTIMASI_WAKTU\", \"WAKTU\"],\n  [\"KRITERIA_PENILAIAN\", \"KRITERIA\"],\n  [\"BOBOT_NILAI\", \"BOBOT\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldName, prefix] of map) {\n  const lines = [];\n  for (let i = 1; i <= 16; i++) {\n    lines.push(`{${prefix}_${i}}`);\n  }\n  // A literal newline character joins the 16 placeholders, all inside a\n  // single <w:t> (NOT a paragraph break / manual line break), exactly as\n  // in the target template.\n  const newText = lines.join(\"\\n\");\n\n  const results = body.search(`{${oldName}}`, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Rebuild the run via a minimal Flat-OPC OOXML fragment so we can\n    // place a raw \"\\n\" inside <w:t> (Range.insertText treats \"\\n\"/\"\\r\" as\n    // a paragraph break, which is not what the template expects here) and\n    // so that the original run formatting (rFonts/sz/szCs/lang) survives\n    // the replace untouched, just like the rest of the run is untouched\n    // in the diff.\n    const flatOpcXml =\n      '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      \"<pkg:xmlData>\" +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      \"<w:body><w:p>\" +\n      '<w:r w:rsidRPr=\"00325ECE\">' +\n      \"<w:rPr>\" +\n      '<w:rFonts w:eastAsia=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n      '<w:sz w:val=\"20\"/>' +\n      '<w:szCs w:val=\"20\"/>' +\n      '<w:lang w:val=\"id-ID\" w:eastAsia=\"id-ID\"/>' +\n      \"</w:rPr>\" +\n      `<w:t>${newText}</w:t>` +\n      \"</w:r>\" +\n      \"</w:p></w:body></w:document>\" +\n      \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n    results.items[0].insertOoxml(flatOpcXml, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Expand each single weekly placeholder in the Rencana Pembelajaran table\n# into 16 per-row placeholders (one per week), each on its own line inside\n# the same run/paragraph (a literal line-feed inside the <w:t>, matching\n# the template's updated per-row placeholder scheme).\n\n$d = $word.ActiveDocument\n\n# old placeholder name -> new per-week prefix\n$map = @{\n    \"MINGGU_KE\"            = \"MINGGU\"\n    \"SUB_CPMK_MINGGUAN\"    = \"SUB_CPMK\"\n    \"INDIKATOR_MINGGUAN\"   = \"INDIKATOR\"\n    \"TOPIK_MINGGUAN\"       = \"TOPIK\"\n    \"METODE_PEMBELAJARAN\"  = \"METODE\"\n    \"ESTIMASI_WAKTU\"       = \"WAKTU\"\n    \"KRITERIA_PENILAIAN\"   = \"KRITERIA\"\n    \"BOBOT_NILAI\"          = \"BOBOT\"\n}\n\n$newline = [char]10\n\nforeach ($old in $map.Keys) {\n    $prefix = $map[$old]\n\n    $lines = @()\n    for ($i = 1; $i -le 16; $i++) {\n        $lines += \"{\" + $prefix + \"_\" + $i + \"}\"\n    }\n    $replacement = [string]::Join($newline, $lines)\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = \"{\" + $old + \"}\"\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replacement\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null\n}\n"}
